# "Se finaliza el desarrollo" -- relabel the last header column and move
# the active selection, per the recorded edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: column H used to read "Empresa Id"; it is now "Entidad Id".
# (Column E stays "EmailSSO" -- unchanged in content, even though the
# underlying shared-string slot it points at is renumbered.)
$ws.Range("H1").Value = "Entidad Id"

# The sheet's active cell/selection moved from F9 to G2.
$ws.Range("G2").Select()
